$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("B2").Value = 42
$ws.Range("F2").Value = 92
$ws.Range("I2").Value = 114
$ws.Range("B3").Value = 75
$ws.Range("D3").Value = 134
$ws.Range("E3").Value = 146
$ws.Range("F3").Value = 137
$ws.Range("H3").Value = 154
$ws.Range("B6").Value = 375
$ws.Range("D6").Value = 414
$ws.Range("F6").Value = 537
$ws.Range("H6").Value = 440
$ws.Range("B7").Value = 503
$ws.Range("D7").Value = 646
$ws.Range("E7").Value = 696
$ws.Range("F7").Value = 776
$ws.Range("H7").Value = 715
$ws.Range("I7").Value = 830

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 4
$ws.Range("F5").Value = 15
$ws.Range("F7").Value = 12
$ws.Range("H10").Value = 6
$ws.Range("H18").Value = 2
$ws.Range("B19").Value = 7
$ws.Range("F19").Value = 24
$ws.Range("F21").Value = 12
$ws.Range("F29").Value = 14
$ws.Range("B30").Value = 6
$ws.Range("B32").Value = 16
$ws.Range("E53").Value = 84
$ws.Range("D65").Value = 26
$ws.Range("F65").Value = 38
$ws.Range("D70").Value = 9
$ws.Range("H70").Value = 15
$ws.Range("F77").Value = 22
$ws.Range("F79").Value = 10
$ws.Range("B80").Value = 16
$ws.Range("H88").Value = 7
$ws.Range("F95").Value = 4
$ws.Range("B98").Value = 503
$ws.Range("D98").Value = 646
$ws.Range("E98").Value = 696
$ws.Range("F98").Value = 776
$ws.Range("H98").Value = 715
$ws.Range("I98").Value = 830

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("F6").Value = 11
$ws.Range("F7").Value = 22

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("F5").Value = 8
$ws.Range("F6").Value = 12

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 6

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("F3").Value = 3
$ws.Range("F7").Value = 12

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("B2").Value = 1
$ws.Range("B7").Value = 16

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("F5").Value = 11
$ws.Range("F6").Value = 15

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("H3").Value = 3
$ws.Range("H5").Value = 7

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("B2").Value = 1
$ws.Range("B5").Value = 16

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("E3").Value = 15
$ws.Range("E7").Value = 84

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("F3").Value = 2
$ws.Range("F6").Value = 10

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("D5").Value = 25
$ws.Range("F5").Value = 31
$ws.Range("D6").Value = 26
$ws.Range("F6").Value = 38

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("F5").Value = 13
$ws.Range("F6").Value = 14

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("H5").Value = 5
$ws.Range("H6").Value = 6

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("B3").Value = 2
$ws.Range("F5").Value = 17
$ws.Range("B6").Value = 7
$ws.Range("F6").Value = 24

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("G2").Value = 2
$ws.Range("G6").Value = 4

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("D3").Value = 3
$ws.Range("H3").Value = 4
$ws.Range("D5").Value = 9
$ws.Range("H5").Value = 15

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("E2").Value = 1
$ws.Range("E5").Value = 4
